# Gate_Closure_Trigger.xlsx update — add "Houma Navigation Canal" gate row.
#
# The row for "WestPoint" / "MS River at West Point a la Hache" (and the
# other MS-River rows that follow it) used to start at worksheet row 23.
# A new row for the Houma Navigation Canal gate is inserted immediately
# above that block (new row 23), pushing the existing rows 23-29 down to
# 24-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 23 ("WestPoint" row), which
# shifts the old rows 23-29 down to 24-30.
[void]$ws.Rows(23).Insert()

# Populate the new row with the Houma Navigation Canal gate entry.
$ws.Range("A23").Value = "HoumaNavCanal"
$ws.Range("B23").Value = "Houma Navigation Canal"
$ws.Range("C23").Value = 0

# Match the formatting of the rows directly below (the "MS River at ..."
# block), which carry a distinct number-format/font style on column A and C.
[void]$ws.Range("A24").Copy()
[void]$ws.Range("A23").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
# Re-assert the text value after the format-only paste.
$ws.Range("A23").Value = "HoumaNavCanal"

# Update the selected cell to match the post-edit selection.
[void]$ws.Range("B23").Select()
